$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 in I1 and IF in J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (style) from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-40
$values = @(
    @(5, 5),
    @(9, 9),
    @(7, 7),
    @(5, 6),
    @(7, 7),
    @(9, 9),
    @(3, 3),
    @(6, 7),
    @(9, 9),
    @(4, 5),
    @(7, 7),
    @(6, 6),
    @(3, 4),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(2, 2),
    @(1, 1),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(7, 8),
    @(6, 6),
    @(8, 8),
    @(5, 6),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(9, 9),
    @(6, 6),
    @(5, 5),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
